$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G width (so the sheet's used range / dimension grows to F, and
# layout matches the author's added column). The host quantizes ColumnWidth
# onto a 1/6-character pixel grid, so 14.666666666666666 is the closest
# input that reproduces the authored stored width of 15.42578125 (-> 15.5).
$ws.Columns.Item(7).ColumnWidth = 14.666666666666666

# Add the new header "Real Effort" in F1, copying the header formatting
# (style) from the existing header cell E1 so it matches the other headers.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Real Effort"

# Move / record the active selection as it was left in the authored file.
$ws.Range("G1").Select()
